$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 867.1739
$ws.Range("J28").Value = 1526.4546
$ws.Range("L28").Value = 1526.4546
$ws.Range("N28").Value = -2496.4546

$ws.Range("H100").Value = 1431.1538
$ws.Range("I100").Value = 1219.9
$ws.Range("J100").Value = 2135.3333
$ws.Range("K100").Value = 1219.9
$ws.Range("L100").Value = 2135.3333
$ws.Range("M100").Value = -678.9000000000001
$ws.Range("N100").Value = -3217.3333

$ws.Range("H107").Value = 779.4706
$ws.Range("I107").Value = 778.25
$ws.Range("J107").Value = 799
$ws.Range("K107").Value = 778.25
$ws.Range("L107").Value = 799
$ws.Range("M107").Value = 1141.75
$ws.Range("N107").Value = -4639

$ws.Range("H111").Value = 2805.9524
$ws.Range("I111").Value = 2445.6155
$ws.Range("J111").Value = 3391.5
$ws.Range("K111").Value = 7336.8465
$ws.Range("L111").Value = 10174.5
$ws.Range("M111").Value = -4269.8465
$ws.Range("N111").Value = -16308.5

$ws.Range("H137").Value = 2438.2896
$ws.Range("I137").Value = 2075.8572
$ws.Range("J137").Value = 6666.6665
$ws.Range("K137").Value = 6227.571599999999
$ws.Range("L137").Value = 19999.9995
$ws.Range("M137").Value = -3677.571599999999
$ws.Range("N137").Value = -25099.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6951.8423
$ws.Range("I74").Value = 4950.375
$ws.Range("K74").Value = 4950.375
$ws.Range("M74").Value = -4076.375

$ws.Range("H77").Value = 6951.8423
$ws.Range("I77").Value = 4950.375
$ws.Range("K77").Value = 24751.875
$ws.Range("M77").Value = -20383.875

$ws.Range("H88").Value = 3192.6667
$ws.Range("I88").Value = 2353
$ws.Range("J88").Value = 3612.5
$ws.Range("K88").Value = 2353
$ws.Range("L88").Value = 3612.5
$ws.Range("M88").Value = -1947
$ws.Range("N88").Value = -4424.5

$ws.Range("H91").Value = 3192.6667
$ws.Range("I91").Value = 2353
$ws.Range("J91").Value = 3612.5
$ws.Range("K91").Value = 2353
$ws.Range("L91").Value = 3612.5
$ws.Range("M91").Value = -949
$ws.Range("N91").Value = -6420.5

$ws.Range("H110").Value = 1322.6522
$ws.Range("I110").Value = 1279.9445
$ws.Range("J110").Value = 1476.4
$ws.Range("K110").Value = 1279.9445
$ws.Range("L110").Value = 1476.4
$ws.Range("M110").Value = 765.0554999999999
$ws.Range("N110").Value = -5566.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1500
$ws.Range("I8").Value = 1500
$ws.Range("K8").Value = 1500
$ws.Range("M8").Value = -1360

$ws.Range("H86").Value = 1811.5714
$ws.Range("I86").Value = 1641
$ws.Range("J86").Value = 9999
$ws.Range("K86").Value = 1641
$ws.Range("L86").Value = 9999
$ws.Range("M86").Value = -518
$ws.Range("N86").Value = -12245

$ws.Range("H89").Value = 1811.5714
$ws.Range("I89").Value = 1641
$ws.Range("J89").Value = 9999
$ws.Range("K89").Value = 8205
$ws.Range("L89").Value = 49995
$ws.Range("M89").Value = -2589
$ws.Range("N89").Value = -61227

$ws.Range("H107").Value = 1087.1428
$ws.Range("I107").Value = 727.7273
$ws.Range("J107").Value = 1482.5
$ws.Range("K107").Value = 727.7273
$ws.Range("L107").Value = 1482.5
$ws.Range("M107").Value = 1192.2727
$ws.Range("N107").Value = -5322.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1154.75
$ws.Range("I16").Value = 858.5714
$ws.Range("J16").Value = 1385.1111
$ws.Range("K16").Value = 858.5714
$ws.Range("L16").Value = 1385.1111
$ws.Range("M16").Value = -571.5714
$ws.Range("N16").Value = -1959.1111

$ws.Range("H31").Value = 4723.077
$ws.Range("I31").Value = 4875.778
$ws.Range("J31").Value = 4379.5
$ws.Range("K31").Value = 4875.778
$ws.Range("L31").Value = 4379.5
$ws.Range("M31").Value = -4580.778
$ws.Range("N31").Value = -4969.5

$ws.Range("H34").Value = 4723.077
$ws.Range("I34").Value = 4875.778
$ws.Range("J34").Value = 4379.5
$ws.Range("K34").Value = 4875.778
$ws.Range("L34").Value = 4379.5
$ws.Range("M34").Value = -4673.778
$ws.Range("N34").Value = -4783.5

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H113").Value = 1154.75
$ws.Range("I113").Value = 858.5714
$ws.Range("J113").Value = 1385.1111
$ws.Range("K113").Value = 858.5714
$ws.Range("L113").Value = 1385.1111
$ws.Range("M113").Value = 1311.4286
$ws.Range("N113").Value = -5725.1111

$ws.Range("H134").Value = 2865.38
$ws.Range("I134").Value = 1771.4231
$ws.Range("J134").Value = 4050.5
$ws.Range("K134").Value = 5314.2693
$ws.Range("L134").Value = 12151.5
$ws.Range("M134").Value = -2779.2693
$ws.Range("N134").Value = -17221.5

$ws.Range("H135").Value = 27187.125
$ws.Range("J135").Value = 27187.125
$ws.Range("L135").Value = 27187.125
$ws.Range("N135").Value = -37327.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 724.1429000000001
$ws.Range("I68").Value = 689.6667
$ws.Range("J68").Value = 750
$ws.Range("K68").Value = 2069.0001
$ws.Range("L68").Value = 2250
$ws.Range("M68").Value = -1258.0001
$ws.Range("N68").Value = -3872

$ws.Range("H71").Value = 724.1429000000001
$ws.Range("I71").Value = 689.6667
$ws.Range("J71").Value = 750
$ws.Range("K71").Value = 6207.0003
$ws.Range("L71").Value = 6750
$ws.Range("M71").Value = -2151.0003
$ws.Range("N71").Value = -14862

$ws.Range("H107").Value = 948.8788
$ws.Range("I107").Value = 331.9375
$ws.Range("J107").Value = 1529.5294
$ws.Range("K107").Value = 995.8125
$ws.Range("L107").Value = 4588.5882
$ws.Range("M107").Value = 924.1875
$ws.Range("N107").Value = -8428.5882

$ws.Range("H117").Value = 1832.7142
$ws.Range("J117").Value = 2062
$ws.Range("L117").Value = 6186
$ws.Range("N117").Value = -13070

$ws.Range("H137").Value = 24299.4
$ws.Range("I137").Value = 1170
$ws.Range("J137").Value = 37309.688
$ws.Range("K137").Value = 3510
$ws.Range("L137").Value = 111929.064
$ws.Range("M137").Value = 1590
$ws.Range("N137").Value = -122129.064

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 394.83334
$ws.Range("I107").Value = 180.2
$ws.Range("J107").Value = 663.125
$ws.Range("K107").Value = 180.2
$ws.Range("L107").Value = 663.125
$ws.Range("M107").Value = 1739.8
$ws.Range("N107").Value = -4503.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3948.0378
$ws.Range("I136").Value = 1940.6
$ws.Range("J136").Value = 7851.3887
$ws.Range("K136").Value = 5821.799999999999
$ws.Range("L136").Value = 23554.1661
$ws.Range("M136").Value = -3271.799999999999
$ws.Range("N136").Value = -28654.1661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 15690.5
$ws.Range("J123").Value = 15690.5
$ws.Range("L123").Value = 15690.5
$ws.Range("N123").Value = -25490.5
